$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: new header labels for the UISettings block ---
$ws.Range("G10").Value = " unlockedSkinPowerAsInfoBox"
$ws.Range("H10").Value = "showContinueButtonInUnlockedSkin"
$ws.Range("I10").Value = "initialMapCountdownTriggeredByPlayer"

# Copy the style used by the other header cells in that row (D10:F10) onto the new ones
$ws.Range("D10:F10").Copy()
$ws.Range("G10:I10").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Row 10 grew taller to fit the new wrapped header text
$ws.Rows.Item(10).RowHeight = 186

# --- Row 11: default boolean values (false) for the new settings ---
$ws.Range("G11").Value = $false
$ws.Range("H11").Value = $false
$ws.Range("I11").Value = $false

# Copy the style used by the other value cells in that row (D11:F11) onto the new ones
$ws.Range("D11:F11").Copy()
$ws.Range("G11:I11").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Update the active view/selection to match where the edit was made ---
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("F11:I11").Select()
